$d = $word.ActiveDocument

# Build the replacement body content as raw OOXML (package-wrapped so
# InsertXML can apply it to word/document.xml) covering all eight new
# paragraphs, including the spell-check proofErr markers around the two
# surnames that Word's proofer would have flagged.
$bodyXml = @'
<w:p>
<w:r><w:t>Project Lead: Ryan Darrow</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Created a proper timeline for everyone to follow for the group to meet the deadline</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Sent email updates with everyone&#8217;s progress and goals</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Set up time and organized zoom meetings with the group</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">Technical Lead: Chris </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Krenz</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">Interface Lead: Jayden </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Raphino</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
<w:r><w:t>Specification Lead: Mark Zhu</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Documentation Lead: Matthew Rhee</w:t></w:r>
</w:p>
'@

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($packageXml)
